# Applies the LOM3238 course-sheet update:
#  - "Integrated Project I" (1 credit-aula / 135h / activated 2012) becomes a
#    general "Integrated Project" research-initiation course (4 credits-aula /
#    180h / activated 2023), with new English objective/short-syllabus/syllabus
#    text and reworded evaluation / recovery / bibliography notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Name (English) ---
$ws.Range("B4").Value = "Integrated Project"
$ws.Range("C4").Value = "Integrated Project"

# --- Row 5: Creditos-aula (numeric-looking text -> force text so it keeps
#     being stored as a shared string, not auto-converted to a number) ---
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

# --- Row 7: Carga horaria ---
$ws.Range("B7").Value = "180 h"
$ws.Range("C7").Value = "180 h"

# --- Row 8: Ativacao (date-looking text -> force text so it is not
#     auto-converted into a date serial number) ---
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

# --- Row 11: Objectives (new content, cells were previously empty) ---
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "Introduce students to the principles and methodology of scientific research."
$ws.Range("C11").Value = "Introduce students to the principles and methodology of scientific research."

# --- Row 13: Programa resumido date mirror (same value as row 8) ---
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

# --- Row 14: Short syllabus (new content, cells were previously empty) ---
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").Value = "Initiation into a research project under the guidance of a professor."
$ws.Range("C14").Value = "Initiation into a research project under the guidance of a professor."

# --- Row 16: Syllabus (new content, cells were previously empty) ---
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = "Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course."
$ws.Range("C16").Value = "Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course."

# --- Row 19: Criterio ---
$ws.Range("B19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."
$ws.Range("C19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."

# --- Row 20: Norma de recuperacao ---
$ws.Range("B20").Value = "Nota de avaliação do projeto e demais documentos."
$ws.Range("C20").Value = "Nota de avaliação do projeto e demais documentos."

# --- Row 21: Bibliografia ---
$ws.Range("B21").Value = "Devido às características do curso, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características do curso, não será oferecida recuperação."

$excel.CutCopyMode = $false
